# Update "paises" (countries) data sheet with the latest scrape (15 Jul 2020, 11:35).
# Several countries' totals changed enough to shuffle their rank within the
# (descending, by "Casos totales") list, so a handful of rows swap which
# country they show even though most rows keep their own country in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Julio de 2020 a las 11:35"

# --- Straightforward numeric refreshes (country keeps its row / rank) ---

# Row 4: Estados Unidos
Set-CountryRow 4 @("Estados Unidos", 3545692, 615, 1600910, 1805627, 0, 12, 139155)

# Row 20: Banglades
Set-CountryRow 20 @("Banglades", 193590, 3533, 105523, 85610, 0, 33, 2457)

# Row 29: Indonesia
Set-CountryRow 29 @("Indonesia", 80094, 1522, 39050, 37247, 0, 87, 3797)

# Row 46: Israel
Set-CountryRow 46 @("Israel", 42813, 453, 19734, 22704, 0, 4, 375)

# Row 47: Polonia
Set-CountryRow 47 @("Polonia", 38721, 264, 28492, 8635, 0, 6, 1594)

# Row 62: Austria
Set-CountryRow 62 @("Austria", 19154, 133, 17175, 1269, 0, 1, 710)

# Row 79: Malasia
Set-CountryRow 79 @("Malasia", 8734, 5, 8526, 86, 0, 0, 122)

# Row 121: Lituania
Set-CountryRow 121 @("Lituania", 1882, 7, 1582, 221, 0, 0, 79)

# Row 122: Eslovenia
Set-CountryRow 122 @("Eslovenia", 1878, 19, 1501, 266, 0, 0, 111)

# Row 126: Hong Kong
Set-CountryRow 126 @("Hong Kong", 1589, 19, 1241, 340, 0, 0, 8)

# Row 149: Surinam (only casos activos / recuperados shift by 1)
Set-CountryRow 149 @("Surinam", 801, 0, 544, 239, 0, 0, 18)

# --- Re-ranked rows: El Salvador's new total (10645) overtakes Australia (10487)
#     and Sudan (10417), so those three rows shuffle down while keeping their own
#     (otherwise unchanged) data ---

# Row 74: now El Salvador (updated data, moved up from row 76)
Set-CountryRow 74 @("El Salvador", 10645, 342, 6120, 4239, 0, 8, 286)

# Row 75: now Australia (its data is unchanged, just dropped one rank)
Set-CountryRow 75 @("Australia", 10487, 237, 7928, 2448, 0, 3, 111)

# Row 76: now Sudan (its data is unchanged, just dropped one rank)
Set-CountryRow 76 @("Sudan", 10417, 0, 5579, 4179, 0, 0, 659)

# --- Islas Malvinas / Groenlandia swap (identical stats, labels only swap) ---

# Row 209: now Islas Malvinas
Set-CountryRow 209 @("Islas Malvinas", 13, 0, 13, 0, 0, 0, 0)

# Row 210: now Groenlandia
Set-CountryRow 210 @("Groenlandia", 13, 0, 13, 0, 0, 0, 0)
